$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting existing rows 6-32 down to 7-33.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new weekly record.
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44552
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112031
$ws.Range("G6").Value = "Poroto verde"
$ws.Range("H6").Value = "Magnum"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 34000
$ws.Range("L6").Value = 36000
$ws.Range("M6").Value = 35000
$ws.Range("N6").Value = "$/malla 25 kilos"
$ws.Range("O6").Value = "Región de O'Higgins"
$ws.Range("P6").Value = 1400
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
